$d = $word.ActiveDocument

$replacements = @(
    @{old = "108÷4="; new = "492÷3="},
    @{old = "552÷3="; new = "208÷3="},
    @{old = "452÷9="; new = "993÷3="},
    @{old = "122÷4="; new = "944÷7="},
    @{old = "729÷3="; new = "808÷2="},
    @{old = "506÷3="; new = "689÷9="},
    @{old = "393÷2="; new = "756÷3="},
    @{old = "894÷5="; new = "532÷5="},
    @{old = "898÷6="; new = "442÷2="},
    @{old = "679÷3="; new = "683÷6="},
    @{old = "154÷3="; new = "670÷8="},
    @{old = "951÷7="; new = "360÷7="},
    @{old = "758÷3="; new = "599÷7="},
    @{old = "749÷4="; new = "688÷2="},
    @{old = "357÷2="; new = "284÷8="},
    @{old = "843÷9="; new = "404÷2="},
    @{old = "390÷6="; new = "259÷8="},
    @{old = "380÷9="; new = "978÷8="},
    @{old = "391÷5="; new = "731÷3="},
    @{old = "482÷3="; new = "168÷4="},
    @{old = "618÷4="; new = "877÷2="},
    @{old = "439÷5="; new = "292÷2="},
    @{old = "762÷8="; new = "236÷6="},
    @{old = "619÷7="; new = "720÷5="},
    @{old = "844÷7="; new = "445÷5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
